# Bond dates update: the "as of" date used to compute day-counts moved
# forward by one day (2023-10-16 -> 2023-10-17). Column G ("Dni od
# poprzedniej wypłaty") is the number of days since the previous payout
# date (column F), so it increases by 1. Column I ("Dni do następnej
# wypłaty") is the number of days until the next payout date (column H),
# so it decreases by 1. Update every data row (2..262) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 262

for ($row = 2; $row -le $lastRow; $row++) {
    $fVal = $ws.Cells.Item($row, 6).Value()   # column F
    $gVal = $ws.Cells.Item($row, 7).Value()   # column G
    $hVal = $ws.Cells.Item($row, 8).Value()   # column H
    $iVal = $ws.Cells.Item($row, 9).Value()   # column I

    if ($fVal -ne $null -and $gVal -ne $null) {
        $ws.Cells.Item($row, 7).Value = $gVal + 1
    }

    if ($hVal -ne $null -and $iVal -ne $null) {
        $ws.Cells.Item($row, 9).Value = $iVal - 1
    }
}
